# regen save_data to use K instead of Strike#, regen std/mean, calc and write s_vals
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated K column (G) values, recomputed from raw pitch data (K = strikeouts)
# replacing the previous erroneous "Strike#" derived values.
$newK = @{
    2  = 2
    3  = 3
    4  = 3
    5  = 2
    6  = 4
    7  = 3
    8  = 5
    9  = 2
    10 = 2
    11 = 4
    12 = 3
    13 = 4
    14 = 8
    15 = 6
    16 = 5
    17 = 6
    18 = 4
    20 = 2
}

foreach ($row in $newK.Keys) {
    $ws.Range("G$row").Value = $newK[$row]
}
